# [MOSIP-19230] Virtual regcenter and zone-user mapping added for Resident services.
#
# zone_user.xlsx: add a zone-user mapping row for the resident-services
# service account, and normalize the "is_active" flag to a live =TRUE()
# formula (was a hard-coded boolean literal).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (existing globaladmin mapping): make is_active a formula ---
$ws.Range("D2").Formula = "=TRUE()"

# --- Row 3 (new): zone-user mapping for the resident services service account ---
$ws.Range("A3").Value = "eng"
$ws.Range("B3").Value = "MOR"
$ws.Range("C3").Value = "Service-account-mosip-resident-client"
$ws.Range("D3").Formula = "=TRUE()"
$ws.Range("D3").NumberFormat = $ws.Range("D2").NumberFormat

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 9.42
$ws.Columns.Item(2).ColumnWidth = 9.92
$ws.Columns.Item(3).ColumnWidth = 31.42
$ws.Columns.Item(4).ColumnWidth = 8.09

# --- Default column width for the sheet ---
$ws.StandardWidth = 11.53515625

# --- Selection moves to C7 ---
$null = $ws.Range("C7").Select()
